$wb = $excel.ActiveWorkbook

# --- Update "Date of Analysis" on the "Project Information" sheet ---
$wsProj = $wb.Worksheets.Item("Project Information")
$wsProj.Range("B12").Value = "2022-04-22 12:42:31.702336"

# --- Reorder the readme table columns: (sheet_name, Date, JobNo, Author) -> (Author, JobNo, sheet_name, Date) ---
$wsReadme = $wb.Worksheets.Item("readme")

# Capture the current data (old column order: B=sheet_name, C=Date, D=JobNo, E=Author) before overwriting.
$sheetNames = @{}
for ($r = 2; $r -le 12; $r++) {
    $sheetNames[$r] = $wsReadme.Range("B$r").Text
}

# New header order
$wsReadme.Range("B1").Value = "Author"
$wsReadme.Range("C1").Value = "JobNo"
$wsReadme.Range("D1").Value = "sheet_name"
$wsReadme.Range("E1").Value = "Date"

# New data rows: Author, JobNo, sheet_name, Date (Date refreshed to the new run date)
# Force the "Date" column to be stored as text (it holds a plain numeric-looking
# string like "20220422", same as the original "20220325" cells).
$wsReadme.Range("E2:E12").NumberFormat = "@"
for ($r = 2; $r -le 12; $r++) {
    $wsReadme.Range("B$r").Value = "jovyan"
    $wsReadme.Range("C$r").Value = "/c/e"
    $wsReadme.Range("D$r").Value = $sheetNames[$r]
    $wsReadme.Range("E$r").Value = "20220422"
}

# Keep the underlying ListObject's column metadata in sync with the new header order.
$lo = $wsReadme.ListObjects.Item(1)
$lo.ListColumns.Item(2).Name = "Author"
$lo.ListColumns.Item(3).Name = "JobNo"
$lo.ListColumns.Item(4).Name = "sheet_name"
$lo.ListColumns.Item(5).Name = "Date"
